# Add a new "2022" column (S) to the electric-intensity table, shifting
# the used range from A1:R11 to A1:S11, and refresh several previously
# estimated values in columns P:R with updated figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Bring formatting for the new column S over from column R --------------
# (column R holds the most recently added year, so its per-row formatting,
# incl. borders/number format, is exactly what the new column S needs)
$ws.Range("R3").Copy() | Out-Null
$ws.Range("S3").PasteSpecial(-4122) | Out-Null

$ws.Range("R4").Copy() | Out-Null
$ws.Range("S4").PasteSpecial(-4122) | Out-Null

$ws.Range("R5:R8").Copy() | Out-Null
$ws.Range("S5:S8").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- New "2022" column header ----------------------------------------------
$ws.Range("S4").Value = 2022

# --- Updated figures for existing years (columns P:R) -----------------------
$ws.Range("P5").Value = 23.111083656771282
$ws.Range("Q5").Value = 24.08077930418019
$ws.Range("R5").Value = 19.336931533747723

$ws.Range("P6").Value = 14.322631450320875
$ws.Range("Q6").Value = 13.073459110725862
$ws.Range("R6").Value = 10.464141365743002

$ws.Range("P7").Value = 23.612622725489956

# --- New figures for 2022 (column S) ----------------------------------------
$ws.Range("S5").Value = 13.600365850576139
$ws.Range("S6").Value = 9.2742414863791556
$ws.Range("S7").Value = 17.303523954725925
$ws.Range("S8").Value = 205.5

# --- Update the saved selection to reflect the new rightmost column --------
$ws.Range("Q15").Select() | Out-Null
